# fix latency units in report sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append units to Read Latency columns (I=min, J=max, K=average) for data rows 3-23.
# Row 7 uses microseconds (usec), all other rows use milliseconds (msec).
for ($row = 3; $row -le 23; $row++) {
    if ($row -eq 7) {
        $unit = "usec"
    } else {
        $unit = "msec"
    }

    foreach ($col in "I", "J", "K") {
        $cell = $ws.Range("$col$row")
        $cell.Value = "$($cell.Value2) $unit"
    }
}
